$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.348.53"
$ws.Range("E2").Value = "  +1.74%  "

$ws.Range("D3").Value = "3.272.68"
$ws.Range("E3").Value = "  +5.32%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "'244.60"
$ws.Range("E5").Value = "  +2.86%  "

$ws.Range("D6").Value = "'614.17"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'1.11"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("D8").Value = "'0.379"
$ws.Range("E8").Value = "  -3.29%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "3.274.87"
$ws.Range("E10").Value = "  +5.50%  "

$ws.Range("D11").Value = "'0.775"
$ws.Range("E11").Value = "  -7.31%  "

$ws.Range("D12").Value = "'0.198"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "96.423.96"
$ws.Range("E13").Value = "  +2.31%  "

$ws.Range("D14").Value = "'0.0000242"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.899.02"
$ws.Range("E15").Value = "  +5.90%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'34.64"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "'5.45"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("D18").Value = "3.307.47"
$ws.Range("E18").Value = "  +6.52%  "

$ws.Range("D19").Value = "'3.51"
$ws.Range("E19").Value = "  -4.59%  "

$ws.Range("D20").Value = "'14.76"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").Value = "'479.84"
$ws.Range("E21").Value = "  +6.37%  "

$ws.Range("D22").Value = "'5.74"
$ws.Range("E22").Value = "  -3.25%  "

$ws.Range("D23").Value = "'0.0000201"
$ws.Range("E23").Value = "  +1.43%  "

$ws.Range("D24").Value = "'9.10"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").Value = "'5.56"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").Value = "'87.39"
$ws.Range("E26").Value = "  +1.63%  "

$ws.Range("D27").Value = "'11.86"
$ws.Range("E27").Value = "  -2.58%  "

$ws.Range("D28").Value = "3.537.04"
$ws.Range("E28").Value = "  +7.71%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("D30").Value = "'0.179"
$ws.Range("E30").Value = "  -2.14%  "

$ws.Range("D31").Value = "'0.237"
$ws.Range("E31").Value = "  -6.93%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").Value = "'0.119"
$ws.Range("E33").Value = "  -4.62%  "

$ws.Range("D34").Value = "'9.15"
$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("D35").Value = "'26.87"
$ws.Range("E35").Value = "  +3.37%  "

$ws.Range("D36").Value = "'7.26"
$ws.Range("E36").Value = "  -7.86%  "

$ws.Range("D37").Value = "'0.149"
$ws.Range("E37").Value = "  -7.30%  "

$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.91"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'24.67"
$ws.Range("E39").Value = "  +2.80%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'488.36"
$ws.Range("E40").Value = "  +3.20%  "

$ws.Range("D41").Value = "'0.444"
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("D42").Value = "'1.25"
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.19"
$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.44"
$ws.Range("E45").Value = "  -7.21%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'0.771"
$ws.Range("E46").Value = "  +11.47%  "

$ws.Range("D47").Value = "'160.85"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").Value = "'1.89"
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'45.07"
$ws.Range("E49").Value = "  +2.86%  "

$ws.Range("D50").Value = "'4.45"
$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "'1.33"
$ws.Range("E51").Value = "  +1.55%  "
